$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = "車貸"
$ws.Range("D2").Value = "台新銀行三重分行新北市三重區正義北路"
$ws.Range("F2").Value = "99年08月27日"
$ws.Range("G2").Value = "車貸"

# Row 3
$ws.Range("B3").Value = "貸款"
$ws.Range("D3").Value = "土地銀行台東分行臺東縣台東市中華路"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "1848023"
$ws.Range("F3").Value = "100年07月01曰"
$ws.Range("G3").Value = "信貸"

# Row 4
$ws.Range("B4").Value = "—般貸款"
$ws.Range("D4").Value = "台東縣都蘭農會臺東縣東河鄉都蘭村都蘭"
$ws.Range("F4").Value = "96年03月13曰"
$ws.Range("G4").Value = "—般貸款"
